$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so number-like strings
# (e.g. "0.996", "47.60") are not auto-converted to numeric values by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '67.389.91'
$ws.Range("E2").Value = '  +1.33%  '

# Row 3
$ws.Range("D3").Value = '3.370.28'
$ws.Range("E3").Value = '  +0.50%  '

# Row 4
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").Value = '  -0.31%  '

# Row 5
$ws.Range("D5").Value = '590.95'
$ws.Range("E5").Value = '  +5.72%  '

# Row 6
$ws.Range("D6").Value = '188.44'
$ws.Range("E6").Value = '  -0.90%  '

# Row 7
$ws.Range("E7").Value = '  +0.12%  '

# Row 8
$ws.Range("D8").Value = '0.599'
$ws.Range("E8").Value = '  +2.54%  '

# Row 9
$ws.Range("E9").Value = '  +0.89%  '

# Row 10
$ws.Range("D10").Value = '0.588'
$ws.Range("E10").Value = '  +0.17%  '

# Row 11
$ws.Range("D11").Value = '47.60'
$ws.Range("E11").Value = '  +0.71%  '

# Row 12
$ws.Range("D12").Value = '0.0000275'
$ws.Range("E12").Value = '  +1.61%  '

# Row 13
$ws.Range("D13").Value = '636.88'
$ws.Range("E13").Value = '  +5.13%  '

# Row 14
$ws.Range("D14").Value = '3.903.83'
$ws.Range("E14").Value = '  +1.04%  '

# Row 15
$ws.Range("E15").Value = '  -0.66%  '

# Row 16
$ws.Range("D16").Value = '67.208.58'
$ws.Range("E16").Value = '  +1.16%  '

# Row 17
$ws.Range("E17").Value = '  +0.75%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.363.15'
$ws.Range("E18").Value = '  +1.08%  '

# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '18.05'
$ws.Range("E19").Value = '  +0.09%  '

# Row 20
$ws.Range("D20").Value = '11.21'
$ws.Range("E20").Value = '  +0.75%  '

# Row 21
$ws.Range("D21").Value = '0.912'
$ws.Range("E21").Value = '  +0.70%  '

# Row 22
$ws.Range("D22").Value = '18.08'
$ws.Range("E22").Value = '  -1.81%  '

# Row 23
$ws.Range("D23").Value = '5.11'
$ws.Range("E23").Value = '  +0.35%  '

# Row 24
$ws.Range("D24").Value = '100.44'
$ws.Range("E24").Value = '  -0.40%  '

# Row 25
$ws.Range("E25").Value = '  +1.08%  '

# Row 26
$ws.Range("E26").Value = '  +3.06%  '

# Row 27
$ws.Range("D27").Value = '9.76'
$ws.Range("E27").Value = '  +1.04%  '

# Row 28
$ws.Range("D28").Value = '32.64'
$ws.Range("E28").Value = '  +6.43%  '

# Row 29
$ws.Range("D29").Value = '8.74'
$ws.Range("E29").Value = '  -0.24%  '

# Row 30
$ws.Range("E30").Value = '  +0.51%  '

# Row 31
$ws.Range("D31").Value = '613.88'
$ws.Range("E31").Value = '  +5.45%  '

# Row 32
$ws.Range("D32").Value = '3.86'
$ws.Range("E32").Value = '  -1.04%  '

# Row 33
$ws.Range("D33").Value = '11.20'
$ws.Range("E33").Value = '  +1.03%  '

# Row 34
$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '3.923.34'
$ws.Range("E34").Value = '  +5.71%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.107'
$ws.Range("E35").Value = '  +2.04%  '

# Row 36
$ws.Range("E36").Value = '  +0.06%  '

# Row 37
$ws.Range("D37").Value = '55.96'
$ws.Range("E37").Value = '  -2.24%  '

# Row 38
$ws.Range("D38").Value = '2.83'
$ws.Range("E38").Value = '  +5.57%  '

# Row 39
$ws.Range("D39").Value = '0.132'
$ws.Range("E39").Value = '  +2.47%  '

# Row 40
$ws.Range("D40").Value = '33.86'
$ws.Range("E40").Value = '  -0.15%  '

# Row 41
$ws.Range("D41").Value = '3.26'
$ws.Range("E41").Value = '  -0.09%  '

# Row 42
$ws.Range("D42").Value = '0.0₃0707'
$ws.Range("E42").Value = '  -0.24%  '

# Row 43
$ws.Range("E43").Value = '  +1.28%  '

# Row 44
$ws.Range("E44").Value = '  -0.78%  '

# Row 45
$ws.Range("D45").Value = '0.0424'
$ws.Range("E45").Value = '  +0.82%  '

# Row 46
$ws.Range("E46").Value = '  +0.15%  '

# Row 47
$ws.Range("D47").Value = '2.59'
$ws.Range("E47").Value = '  +0.23%  '

# Row 48
$ws.Range("E48").Value = '  +0.90%  '

# Row 49
$ws.Range("D49").Value = '2.87'
$ws.Range("E49").Value = '  -19.32%  '

# Row 50
$ws.Range("E50").Value = '  +8.51%  '

# Row 51
$ws.Range("D51").Value = '129.38'
$ws.Range("E51").Value = '  +4.00%  '
